# regmap table minor fix
# - "Digital Outputs (Write)" row becomes "Digital Outputs (Read/Write)"
# - "Digital Outputs (Set/Clear/Toggle)" rows get a trailing "*"
# - a footnote "* - write only" is added below the Detail table
# - the Detail table gets a proper two-column header ("Detail" spanning
#   I2:J2, bordered/bold/centered like the Summary header, and I3:J3
#   bordered/bold like the Summary table's column headers)
# - column J is widened a bit to fit the new header text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the Digital Outputs descriptions -------------------------
$ws.Range("J17").Value = "Digital Outputs (Read/Write)"
$ws.Range("J18").Value = "Digital Outputs (Set)*"
$ws.Range("J19").Value = "Digital Outputs (Clear)*"
$ws.Range("J20").Value = "Digital Outputs (Toggle)*"

# --- 2. Add the footnote row ---------------------------------------------
$ws.Range("J22").Value = "* - write only"

# --- 3. Widen column J so the longer text still fits ---------------------
$ws.Columns.Item(10).ColumnWidth = 27.85546875

# --- 4. Format + merge the "Detail" header (I2:J2) like the "Summary" one
# (merge first, then paste borders, so Excel doesn't "smart-derive" open
#  merged-cell edges the way it would if we bordered first and merged after)
$ws.Range("I2:J2").Merge()
$ws.Range("C4").Copy()
$ws.Range("I2:J2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$hdr = $ws.Range("I2:J2")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108

# --- 5. Format the "Reg. Address" / "Description" column headers (I3:J3)
#        like the Summary table's "Range/Type/Function" headers
$ws.Range("C3").Copy()
$ws.Range("I3:J3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 6. Cosmetic: move the active selection like in the authored file ----
$ws.Range("I26").Select() | Out-Null
